# Updated cryptos list on Fri Aug 18 23:55:16 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and 1h-volume-change (column E) figures refreshed for
# each coin row. D values are assigned with a leading apostrophe so Excel
# keeps them as literal text (matching the sheet's existing string-typed
# price cells) instead of re-parsing them as numbers and silently dropping
# meaningful trailing zeros (e.g. "4.540" -> 4.54).

$ws.Range("D2").Value = "'26.231.68"
$ws.Range("E2").Value = "  -2.00%  "

$ws.Range("D3").Value = "'1.670.67"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'217.53"
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("D6").Value = "'0.5116"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.2655"
$ws.Range("E8").Value = "  +3.32%  "

$ws.Range("D9").Value = "'0.06376"
$ws.Range("E9").Value = "  +3.31%  "

$ws.Range("D10").Value = "'21.49"
$ws.Range("E10").Value = "  -1.64%  "

$ws.Range("D11").Value = "'0.07388"
$ws.Range("E11").Value = "  +0.64%  "

$ws.Range("D12").Value = "'4.540"
$ws.Range("E12").Value = "  +1.74%  "

$ws.Range("D13").Value = "'1.673.77"
$ws.Range("E13").Value = "  -1.28%  "

$ws.Range("D14").Value = "'0.5815"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("D15").Value = "'0.000008638"
$ws.Range("E15").Value = "  +5.56%  "

$ws.Range("D16").Value = "'64.46"
$ws.Range("E16").Value = "  -1.03%  "

$ws.Range("D17").Value = "'26.201.57"
$ws.Range("E17").Value = "  -2.18%  "

$ws.Range("D18").Value = "'4.936"
$ws.Range("E18").Value = "  -1.77%  "

$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").Value = "'10.87"
$ws.Range("E20").Value = "  +2.10%  "

$ws.Range("D21").Value = "'189.20"
$ws.Range("E21").Value = "  +1.40%  "

$ws.Range("D22").Value = "'6.208"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").Value = "'144.33"
$ws.Range("E24").Value = "  +1.27%  "

$ws.Range("D25").Value = "'7.638"
$ws.Range("E25").Value = "  +2.16%  "

$ws.Range("D26").Value = "'0.1176"
$ws.Range("E26").Value = "  +2.84%  "

$ws.Range("D27").Value = "'15.64"
$ws.Range("E27").Value = "  +3.14%  "

$ws.Range("D28").Value = "'0.05976"
$ws.Range("E28").Value = "  +1.92%  "

$ws.Range("D29").Value = "'1.286"
$ws.Range("E29").Value = "  -3.49%  "

$ws.Range("D30").Value = "'1.324"
$ws.Range("E30").Value = "  -1.52%  "

$ws.Range("D31").Value = "'3.523"
$ws.Range("E31").Value = "  +1.85%  "

$ws.Range("D32").Value = "'3.522"
$ws.Range("E32").Value = "  +2.74%  "

$ws.Range("D33").Value = "'1.643"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").Value = "'1.013"
$ws.Range("E34").Value = "  +2.03%  "

$ws.Range("D35").Value = "'0.6030"
$ws.Range("E35").Value = "  +1.21%  "

$ws.Range("D36").Value = "'2.375"
$ws.Range("E36").Value = "  -1.68%  "

$ws.Range("D37").Value = "'2.656"
$ws.Range("E37").Value = "  -0.40%  "

$ws.Range("E38").Value = "  +4.66%  "

$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("D40").Value = "'1.077.31"
$ws.Range("E40").Value = "  -1.19%  "

$ws.Range("D41").Value = "'0.8695"
$ws.Range("E41").Value = "  +1.60%  "

$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("D43").Value = "'100.30"
$ws.Range("E43").Value = "  +3.01%  "

$ws.Range("D44").Value = "'1.820.06"
$ws.Range("E44").Value = "  -1.18%  "

$ws.Range("D45").Value = "'0.00000000112"
$ws.Range("E45").Value = "  +8.30%  "

$ws.Range("D46").Value = "'56.28"
$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("E47").Value = "  +0.98%  "

$ws.Range("D48").Value = "'8.061"
$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").Value = "'5.875"
$ws.Range("E51").Value = "  +2.36%  "

# The apostrophe prefix marks a cell as quote-prefixed text, which would
# otherwise tag these cells with an extra style. Clear that incidental
# formatting so D2:D51 keep the same (unstyled) look they had before.
$ws.Range("D2:D51").ClearFormats()
